$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap home/away match data between re-ordered row pairs ---
# Row 6
$ws.Cells.Item(6, 6).Value = 'Tuzlaspor'
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 'Corum'
$ws.Cells.Item(6, 9).Value = 3
$ws.Cells.Item(6, 10).Value = 2.16
$ws.Cells.Item(6, 11).Value = '12/08/2023 10:38'
$ws.Cells.Item(6, 12).Value = 2.7
$ws.Cells.Item(6, 13).Value = '13/08/2023 15:54'
$ws.Cells.Item(6, 14).Value = 3.19
$ws.Cells.Item(6, 15).Value = '12/08/2023 10:38'
$ws.Cells.Item(6, 16).Value = 3.24
$ws.Cells.Item(6, 17).Value = '13/08/2023 15:54'
$ws.Cells.Item(6, 18).Value = 3.4
$ws.Cells.Item(6, 19).Value = '12/08/2023 10:38'
$ws.Cells.Item(6, 20).Value = 2.74
$ws.Cells.Item(6, 21).Value = '13/08/2023 15:54'
$ws.Cells.Item(6, 22).Value = 'https://www.betexplorer.com/football/turkey/1-lig/tuzlaspor-corum-fk/GfFXhfNQ/'

# Row 7
$ws.Cells.Item(7, 6).Value = 'Genclerbirligi'
$ws.Cells.Item(7, 7).Value = 2
$ws.Cells.Item(7, 8).Value = 'Umraniyespor'
$ws.Cells.Item(7, 9).Value = 1
$ws.Cells.Item(7, 10).Value = 2.34
$ws.Cells.Item(7, 11).Value = '12/08/2023 10:38'
$ws.Cells.Item(7, 12).Value = 2.93
$ws.Cells.Item(7, 13).Value = '13/08/2023 15:52'
$ws.Cells.Item(7, 14).Value = 3.26
$ws.Cells.Item(7, 15).Value = '12/08/2023 10:38'
$ws.Cells.Item(7, 16).Value = 3.38
$ws.Cells.Item(7, 17).Value = '13/08/2023 15:52'
$ws.Cells.Item(7, 18).Value = 2.97
$ws.Cells.Item(7, 19).Value = '12/08/2023 10:38'
$ws.Cells.Item(7, 20).Value = 2.46
$ws.Cells.Item(7, 21).Value = '13/08/2023 15:52'
$ws.Cells.Item(7, 22).Value = 'https://www.betexplorer.com/football/turkey/1-lig/genclerbirligi-umraniyespor/pAqYXho1/'

# Row 15
$ws.Cells.Item(15, 6).Value = 'Bandirmaspor'
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 'Boluspor'
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 1.86
$ws.Cells.Item(15, 11).Value = '13/08/2023 16:13'
$ws.Cells.Item(15, 12).Value = 1.73
$ws.Cells.Item(15, 13).Value = '20/08/2023 15:53'
$ws.Cells.Item(15, 14).Value = 3.66
$ws.Cells.Item(15, 15).Value = '13/08/2023 16:13'
$ws.Cells.Item(15, 16).Value = 3.82
$ws.Cells.Item(15, 17).Value = '20/08/2023 15:53'
$ws.Cells.Item(15, 18).Value = 4.06
$ws.Cells.Item(15, 19).Value = '13/08/2023 16:13'
$ws.Cells.Item(15, 20).Value = 4.79
$ws.Cells.Item(15, 21).Value = '20/08/2023 15:53'
$ws.Cells.Item(15, 22).Value = 'https://www.betexplorer.com/football/turkey/1-lig/bandirmaspor-boluspor/zLm1jMIP/'

# Row 16
$ws.Cells.Item(16, 6).Value = 'Manisa FK'
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(16, 8).Value = 'Keciorengucu'
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = 1.97
$ws.Cells.Item(16, 11).Value = '13/08/2023 16:13'
$ws.Cells.Item(16, 12).Value = 1.87
$ws.Cells.Item(16, 13).Value = '20/08/2023 15:54'
$ws.Cells.Item(16, 14).Value = 3.63
$ws.Cells.Item(16, 15).Value = '13/08/2023 16:13'
$ws.Cells.Item(16, 16).Value = 3.58
$ws.Cells.Item(16, 17).Value = '20/08/2023 15:59'
$ws.Cells.Item(16, 18).Value = 3.66
$ws.Cells.Item(16, 19).Value = '13/08/2023 16:13'
$ws.Cells.Item(16, 20).Value = 4.27
$ws.Cells.Item(16, 21).Value = '20/08/2023 15:54'
$ws.Cells.Item(16, 22).Value = 'https://www.betexplorer.com/football/turkey/1-lig/manisa-fk-keciorengucu/lSTypKeg/'

# Row 31
$ws.Cells.Item(31, 6).Value = 'Boluspor'
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 'Keciorengucu'
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 10).Value = 2.78
$ws.Cells.Item(31, 11).Value = '27/08/2023 18:42'
$ws.Cells.Item(31, 12).Value = 2.43
$ws.Cells.Item(31, 13).Value = '02/09/2023 18:13'
$ws.Cells.Item(31, 14).Value = 3.3
$ws.Cells.Item(31, 15).Value = '27/08/2023 18:42'
$ws.Cells.Item(31, 16).Value = 3.33
$ws.Cells.Item(31, 17).Value = '02/09/2023 18:07'
$ws.Cells.Item(31, 18).Value = 2.56
$ws.Cells.Item(31, 19).Value = '27/08/2023 18:42'
$ws.Cells.Item(31, 20).Value = 3.01
$ws.Cells.Item(31, 21).Value = '02/09/2023 18:13'
$ws.Cells.Item(31, 22).Value = 'https://www.betexplorer.com/football/turkey/1-lig/boluspor-keciorengucu/WSWSskWM/'

# Row 32
$ws.Cells.Item(32, 6).Value = 'Manisa FK'
$ws.Cells.Item(32, 7).Value = 4
$ws.Cells.Item(32, 8).Value = 'Tuzlaspor'
$ws.Cells.Item(32, 9).Value = 0
$ws.Cells.Item(32, 10).Value = 1.53
$ws.Cells.Item(32, 11).Value = '26/08/2023 18:13'
$ws.Cells.Item(32, 12).Value = 1.47
$ws.Cells.Item(32, 13).Value = '02/09/2023 18:08'
$ws.Cells.Item(32, 14).Value = 4.28
$ws.Cells.Item(32, 15).Value = '26/08/2023 18:13'
$ws.Cells.Item(32, 16).Value = 4.5
$ws.Cells.Item(32, 17).Value = '02/09/2023 18:08'
$ws.Cells.Item(32, 18).Value = 5.83
$ws.Cells.Item(32, 19).Value = '26/08/2023 18:13'
$ws.Cells.Item(32, 20).Value = 6.82
$ws.Cells.Item(32, 21).Value = '02/09/2023 18:08'
$ws.Cells.Item(32, 22).Value = 'https://www.betexplorer.com/football/turkey/1-lig/manisa-fk-tuzlaspor/jX3gbAG3/'

# Row 60
$ws.Cells.Item(60, 6).Value = 'Bodrumspor'
$ws.Cells.Item(60, 7).Value = 2
$ws.Cells.Item(60, 8).Value = 'Corum'
$ws.Cells.Item(60, 9).Value = 1
$ws.Cells.Item(60, 10).Value = 1.84
$ws.Cells.Item(60, 11).Value = '28/09/2023 03:13'
$ws.Cells.Item(60, 12).Value = 1.83
$ws.Cells.Item(60, 13).Value = '01/10/2023 12:28'
$ws.Cells.Item(60, 14).Value = 3.63
$ws.Cells.Item(60, 15).Value = '28/09/2023 03:13'
$ws.Cells.Item(60, 16).Value = 3.45
$ws.Cells.Item(60, 17).Value = '01/10/2023 12:28'
$ws.Cells.Item(60, 18).Value = 4.22
$ws.Cells.Item(60, 19).Value = '28/09/2023 03:13'
$ws.Cells.Item(60, 20).Value = 4.72
$ws.Cells.Item(60, 21).Value = '01/10/2023 12:28'
$ws.Cells.Item(60, 22).Value = 'https://www.betexplorer.com/football/turkey/1-lig/bodrumspor-corum-fk/hWGJL4ER/'

# Row 61
$ws.Cells.Item(61, 6).Value = 'Erzurumspor'
$ws.Cells.Item(61, 7).Value = 1
$ws.Cells.Item(61, 8).Value = 'Boluspor'
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 10).Value = 2.46
$ws.Cells.Item(61, 11).Value = '24/09/2023 17:13'
$ws.Cells.Item(61, 12).Value = 2.62
$ws.Cells.Item(61, 13).Value = '01/10/2023 12:03'
$ws.Cells.Item(61, 14).Value = 3.24
$ws.Cells.Item(61, 15).Value = '24/09/2023 17:13'
$ws.Cells.Item(61, 16).Value = 3.25
$ws.Cells.Item(61, 17).Value = '01/10/2023 12:23'
$ws.Cells.Item(61, 18).Value = 2.95
$ws.Cells.Item(61, 19).Value = '24/09/2023 17:13'
$ws.Cells.Item(61, 20).Value = 2.82
$ws.Cells.Item(61, 21).Value = '01/10/2023 12:23'
$ws.Cells.Item(61, 22).Value = 'https://www.betexplorer.com/football/turkey/1-lig/erzurumspor-fk-boluspor/256gF2ie/'

# Row 79
$ws.Cells.Item(79, 6).Value = 'Umraniyespor'
$ws.Cells.Item(79, 7).Value = 3
$ws.Cells.Item(79, 8).Value = 'Manisa FK'
$ws.Cells.Item(79, 9).Value = 2
$ws.Cells.Item(79, 10).Value = 2.76
$ws.Cells.Item(79, 11).Value = '16/10/2023 01:12'
$ws.Cells.Item(79, 12).Value = 3.09
$ws.Cells.Item(79, 13).Value = '22/10/2023 14:51'
$ws.Cells.Item(79, 14).Value = 3.29
$ws.Cells.Item(79, 15).Value = '16/10/2023 01:12'
$ws.Cells.Item(79, 16).Value = 3.23
$ws.Cells.Item(79, 17).Value = '22/10/2023 14:51'
$ws.Cells.Item(79, 18).Value = 2.58
$ws.Cells.Item(79, 19).Value = '16/10/2023 01:12'
$ws.Cells.Item(79, 20).Value = 2.43
$ws.Cells.Item(79, 21).Value = '22/10/2023 14:51'
$ws.Cells.Item(79, 22).Value = 'https://www.betexplorer.com/football/turkey/1-lig/umraniyespor-manisa-fk/lKZfXzDt/'

# Row 80
$ws.Cells.Item(80, 6).Value = 'Giresunspor'
$ws.Cells.Item(80, 7).Value = 0
$ws.Cells.Item(80, 8).Value = 'Boluspor'
$ws.Cells.Item(80, 9).Value = 1
$ws.Cells.Item(80, 10).Value = 2.93
$ws.Cells.Item(80, 11).Value = '16/10/2023 20:42'
$ws.Cells.Item(80, 12).Value = 3.53
$ws.Cells.Item(80, 13).Value = '22/10/2023 14:59'
$ws.Cells.Item(80, 14).Value = 3.18
$ws.Cells.Item(80, 15).Value = '16/10/2023 20:42'
$ws.Cells.Item(80, 16).Value = 3.4
$ws.Cells.Item(80, 17).Value = '22/10/2023 14:59'
$ws.Cells.Item(80, 18).Value = 2.52
$ws.Cells.Item(80, 19).Value = '16/10/2023 20:42'
$ws.Cells.Item(80, 20).Value = 2.14
$ws.Cells.Item(80, 21).Value = '22/10/2023 14:59'
$ws.Cells.Item(80, 22).Value = 'https://www.betexplorer.com/football/turkey/1-lig/giresunspor-boluspor/QPUJRWtI/'

# Row 89
$ws.Cells.Item(89, 6).Value = 'Adanaspor AS'
$ws.Cells.Item(89, 7).Value = 1
$ws.Cells.Item(89, 8).Value = 'Umraniyespor'
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = 2.55
$ws.Cells.Item(89, 11).Value = '26/10/2023 15:12'
$ws.Cells.Item(89, 12).Value = 3.46
$ws.Cells.Item(89, 13).Value = '29/10/2023 13:52'
$ws.Cells.Item(89, 14).Value = 3.24
$ws.Cells.Item(89, 15).Value = '26/10/2023 15:12'
$ws.Cells.Item(89, 16).Value = 3.39
$ws.Cells.Item(89, 17).Value = '29/10/2023 13:52'
$ws.Cells.Item(89, 18).Value = 2.83
$ws.Cells.Item(89, 19).Value = '26/10/2023 15:12'
$ws.Cells.Item(89, 20).Value = 2.17
$ws.Cells.Item(89, 21).Value = '29/10/2023 13:52'
$ws.Cells.Item(89, 22).Value = 'https://www.betexplorer.com/football/turkey/1-lig/adanaspor-as-umraniyespor/6ynHQjdO/'

# Row 90
$ws.Cells.Item(90, 6).Value = 'Altay'
$ws.Cells.Item(90, 7).Value = 1
$ws.Cells.Item(90, 8).Value = 'Sakaryaspor'
$ws.Cells.Item(90, 9).Value = 2
$ws.Cells.Item(90, 10).Value = 2.9
$ws.Cells.Item(90, 11).Value = '23/10/2023 19:12'
$ws.Cells.Item(90, 12).Value = 4.56
$ws.Cells.Item(90, 13).Value = '29/10/2023 13:52'
$ws.Cells.Item(90, 14).Value = 3.32
$ws.Cells.Item(90, 15).Value = '23/10/2023 19:12'
$ws.Cells.Item(90, 16).Value = 3.66
$ws.Cells.Item(90, 17).Value = '29/10/2023 13:52'
$ws.Cells.Item(90, 18).Value = 2.45
$ws.Cells.Item(90, 19).Value = '23/10/2023 19:12'
$ws.Cells.Item(90, 20).Value = 1.8
$ws.Cells.Item(90, 21).Value = '29/10/2023 13:52'
$ws.Cells.Item(90, 22).Value = 'https://www.betexplorer.com/football/turkey/1-lig/altay-sakaryaspor/4WHotWeB/'

# Row 106
$ws.Cells.Item(106, 6).Value = 'Manisa FK'
$ws.Cells.Item(106, 7).Value = 0
$ws.Cells.Item(106, 8).Value = 'Bandirmaspor'
$ws.Cells.Item(106, 9).Value = 2
$ws.Cells.Item(106, 10).Value = 2.18
$ws.Cells.Item(106, 11).Value = '05/11/2023 11:42'
$ws.Cells.Item(106, 12).Value = 2.54
$ws.Cells.Item(106, 13).Value = '12/11/2023 11:21'
$ws.Cells.Item(106, 14).Value = 3.41
$ws.Cells.Item(106, 15).Value = '05/11/2023 11:42'
$ws.Cells.Item(106, 16).Value = 3.35
$ws.Cells.Item(106, 17).Value = '12/11/2023 11:25'
$ws.Cells.Item(106, 18).Value = 3.31
$ws.Cells.Item(106, 19).Value = '05/11/2023 11:42'
$ws.Cells.Item(106, 20).Value = 2.84
$ws.Cells.Item(106, 21).Value = '12/11/2023 11:25'
$ws.Cells.Item(106, 22).Value = 'https://www.betexplorer.com/football/turkey/1-lig/manisa-fk-bandirmaspor/Eq666S9G/'

# Row 107
$ws.Cells.Item(107, 6).Value = 'Boluspor'
$ws.Cells.Item(107, 7).Value = 2
$ws.Cells.Item(107, 8).Value = 'Sakaryaspor'
$ws.Cells.Item(107, 9).Value = 3
$ws.Cells.Item(107, 10).Value = 2.79
$ws.Cells.Item(107, 11).Value = '05/11/2023 14:12'
$ws.Cells.Item(107, 12).Value = 2.69
$ws.Cells.Item(107, 13).Value = '12/11/2023 11:01'
$ws.Cells.Item(107, 14).Value = 3.15
$ws.Cells.Item(107, 15).Value = '05/11/2023 14:12'
$ws.Cells.Item(107, 16).Value = 3.01
$ws.Cells.Item(107, 17).Value = '12/11/2023 10:41'
$ws.Cells.Item(107, 18).Value = 2.65
$ws.Cells.Item(107, 19).Value = '05/11/2023 14:12'
$ws.Cells.Item(107, 20).Value = 2.95
$ws.Cells.Item(107, 21).Value = '12/11/2023 11:01'
$ws.Cells.Item(107, 22).Value = 'https://www.betexplorer.com/football/turkey/1-lig/boluspor-sakaryaspor/IcM3n7Ip/'

# --- Append two new match rows (119, 120), cloning formatting from row 118 ---
$ws.Range("A118:V118").Copy()
$ws.Range("A119:V120").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 119
$ws.Cells.Item(119, 1).Value = 118
$ws.Cells.Item(119, 2).Value = 'turkey'
$ws.Cells.Item(119, 3).Value = '1-lig'
$ws.Cells.Item(119, 4).Value = '2023-2024'
$ws.Cells.Item(119, 5).Value = 45261.75
$ws.Cells.Item(119, 6).Value = 'Adanaspor AS'
$ws.Cells.Item(119, 7).Value = 0
$ws.Cells.Item(119, 8).Value = 'Corum'
$ws.Cells.Item(119, 9).Value = 3
$ws.Cells.Item(119, 10).Value = 3.06
$ws.Cells.Item(119, 11).Value = '26/11/2023 11:42'
$ws.Cells.Item(119, 12).Value = 3.75
$ws.Cells.Item(119, 13).Value = '01/12/2023 17:58'
$ws.Cells.Item(119, 14).Value = 3.19
$ws.Cells.Item(119, 15).Value = '26/11/2023 11:42'
$ws.Cells.Item(119, 16).Value = 3.53
$ws.Cells.Item(119, 17).Value = '01/12/2023 17:58'
$ws.Cells.Item(119, 18).Value = 2.42
$ws.Cells.Item(119, 19).Value = '26/11/2023 11:42'
$ws.Cells.Item(119, 20).Value = 2.02
$ws.Cells.Item(119, 21).Value = '01/12/2023 17:58'
$ws.Cells.Item(119, 22).Value = 'https://www.betexplorer.com/football/turkey/1-lig/adanaspor-as-corum-fk/nic4nSZs/'

# Row 120
$ws.Cells.Item(120, 1).Value = 119
$ws.Cells.Item(120, 2).Value = 'turkey'
$ws.Cells.Item(120, 3).Value = '1-lig'
$ws.Cells.Item(120, 4).Value = '2023-2024'
$ws.Cells.Item(120, 5).Value = 45261.75
$ws.Cells.Item(120, 6).Value = 'Kocaelispor'
$ws.Cells.Item(120, 7).Value = 0
$ws.Cells.Item(120, 8).Value = 'Umraniyespor'
$ws.Cells.Item(120, 9).Value = 3
$ws.Cells.Item(120, 10).Value = 1.9
$ws.Cells.Item(120, 11).Value = '25/11/2023 18:13'
$ws.Cells.Item(120, 12).Value = 1.76
$ws.Cells.Item(120, 13).Value = '01/12/2023 17:34'
$ws.Cells.Item(120, 14).Value = 3.52
$ws.Cells.Item(120, 15).Value = '25/11/2023 18:13'
$ws.Cells.Item(120, 16).Value = 3.73
$ws.Cells.Item(120, 17).Value = '01/12/2023 17:34'
$ws.Cells.Item(120, 18).Value = 4.06
$ws.Cells.Item(120, 19).Value = '25/11/2023 18:13'
$ws.Cells.Item(120, 20).Value = 4.73
$ws.Cells.Item(120, 21).Value = '01/12/2023 17:49'
$ws.Cells.Item(120, 22).Value = 'https://www.betexplorer.com/football/turkey/1-lig/kocaelispor-umraniyespor/GWlLrpZ6/'
